$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "K" (strikeout) values for rows 2-32 in column G,
# regenerated from source data (per commit message: "use K instead of Strike#").
$kValues = @{
    2  = 2
    3  = 4
    4  = 6
    5  = 3
    6  = 5
    7  = 7
    8  = 5
    9  = 2
    10 = 13
    11 = 6
    12 = 0
    13 = 4
    14 = 3
    15 = 3
    16 = 4
    17 = 1
    18 = 3
    19 = 2
    20 = 9
    21 = 6
    22 = 5
    23 = 4
    24 = 0
    25 = 6
    26 = 4
    27 = 7
    28 = 3
    29 = 4
    30 = 4
    31 = 4
    32 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
